# Auto-generated edit script: updates market-price/profit values in each class leve table
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1719.1936
$ws.Range("I98").Value = 1188.8462
$ws.Range("K98").Value = 1188.8462
$ws.Range("M98").Value = 309.1538

$ws.Range("H107").Value = 1173.6
$ws.Range("I107").Value = 898.75
$ws.Range("J107").Value = 2273
$ws.Range("K107").Value = 898.75
$ws.Range("L107").Value = 2273
$ws.Range("M107").Value = 1021.25
$ws.Range("N107").Value = -6113

$ws.Range("H113").Value = 4000
$ws.Range("J113").Value = 4000
$ws.Range("L113").Value = 4000
$ws.Range("N113").Value = -10508

$ws.Range("H122").Value = 1719.1936
$ws.Range("I122").Value = 1188.8462
$ws.Range("K122").Value = 3566.5386
$ws.Range("M122").Value = -1116.5386

$ws.Range("H137").Value = 1230.45
$ws.Range("I137").Value = 1191.5294
$ws.Range("K137").Value = 3574.5882
$ws.Range("M137").Value = -1024.5882

$ws.Range("H141").Value = 761.6667
$ws.Range("I141").Value = 761.6667
$ws.Range("K141").Value = 2285.0001
$ws.Range("M141").Value = 2894.9999


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3499.75
$ws.Range("I32").Value = 3700.3438
$ws.Range("J32").Value = 1895
$ws.Range("K32").Value = 3700.3438
$ws.Range("L32").Value = 1895
$ws.Range("M32").Value = -3413.3438
$ws.Range("N32").Value = -2469

$ws.Range("H45").Value = 1673.6364
$ws.Range("I45").Value = 1835.1111
$ws.Range("J45").Value = 947
$ws.Range("K45").Value = 1835.1111
$ws.Range("L45").Value = 947
$ws.Range("M45").Value = -1458.1111
$ws.Range("N45").Value = -1701

$ws.Range("H110").Value = 2074.0715
$ws.Range("J110").Value = 2550.5
$ws.Range("L110").Value = 2550.5
$ws.Range("N110").Value = -6640.5

$ws.Range("H122").Value = 2201.5557
$ws.Range("I122").Value = 2176
$ws.Range("J122").Value = 2233.5
$ws.Range("K122").Value = 6528
$ws.Range("L122").Value = 6700.5
$ws.Range("M122").Value = -4078
$ws.Range("N122").Value = -11600.5


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4188.909
$ws.Range("I134").Value = 1155.8148
$ws.Range("J134").Value = 17837.834
$ws.Range("K134").Value = 3467.4444
$ws.Range("L134").Value = 53513.50199999999
$ws.Range("M134").Value = -932.4444000000003
$ws.Range("N134").Value = -58583.50199999999


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 142858420
$ws.Range("I16").Value = 250001100
$ws.Range("J16").Value = 1492.6666
$ws.Range("K16").Value = 250001100
$ws.Range("L16").Value = 1492.6666
$ws.Range("M16").Value = -250000813
$ws.Range("N16").Value = -2066.6666

$ws.Range("H113").Value = 142858420
$ws.Range("I113").Value = 250001100
$ws.Range("J113").Value = 1492.6666
$ws.Range("K113").Value = 250001100
$ws.Range("L113").Value = 1492.6666
$ws.Range("M113").Value = -249998930
$ws.Range("N113").Value = -5832.6666

$ws.Range("H134").Value = 17545400
$ws.Range("I134").Value = 22223872
$ws.Range("K134").Value = 66671616
$ws.Range("M134").Value = -66669081


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1466615.2
$ws.Range("I4").Value = 149923
$ws.Range("J4").Value = 4099999.8
$ws.Range("K4").Value = 449769
$ws.Range("L4").Value = 12299999.4
$ws.Range("M4").Value = -449657
$ws.Range("N4").Value = -12300223.4

$ws.Range("H6").Value = 888.2222
$ws.Range("I6").Value = 133
$ws.Range("J6").Value = 1265.8334
$ws.Range("K6").Value = 399
$ws.Range("L6").Value = 3797.5002
$ws.Range("M6").Value = -286
$ws.Range("N6").Value = -4023.5002

$ws.Range("H7").Value = 450.5
$ws.Range("I7").Value = 533.3333
$ws.Range("J7").Value = 202
$ws.Range("K7").Value = 1599.9999
$ws.Range("L7").Value = 606
$ws.Range("M7").Value = -1487.9999
$ws.Range("N7").Value = -830

$ws.Range("H10").Value = 61.833332
$ws.Range("I10").Value = 34.2
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 102.6
$ws.Range("L10").Value = 600
$ws.Range("M10").Value = 36.39999999999999
$ws.Range("N10").Value = -878

$ws.Range("H11").Value = 193773.5
$ws.Range("I11").Value = 221441.14
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 664323.42
$ws.Range("L11").Value = 300
$ws.Range("M11").Value = -664183.42
$ws.Range("N11").Value = -580

$ws.Range("H37").Value = 200000
$ws.Range("J37").Value = 200000
$ws.Range("L37").Value = 600000
$ws.Range("N37").Value = -600224

$ws.Range("H69").Value = 2313.3125
$ws.Range("I69").Value = 1799.5
$ws.Range("J69").Value = 2386.7144
$ws.Range("K69").Value = 5398.5
$ws.Range("L69").Value = 7160.1432
$ws.Range("M69").Value = -4587.5
$ws.Range("N69").Value = -8782.143199999999

$ws.Range("H72").Value = 2313.3125
$ws.Range("I72").Value = 1799.5
$ws.Range("J72").Value = 2386.7144
$ws.Range("K72").Value = 16195.5
$ws.Range("L72").Value = 21480.4296
$ws.Range("M72").Value = -12139.5
$ws.Range("N72").Value = -29592.4296

$ws.Range("H107").Value = 1138.1428
$ws.Range("J107").Value = 1138.1428
$ws.Range("L107").Value = 3414.4284
$ws.Range("N107").Value = -7254.428400000001

$ws.Range("H122").Value = 956.8570999999999
$ws.Range("I122").Value = 874.5
$ws.Range("K122").Value = 7870.5
$ws.Range("M122").Value = -5420.5

$ws.Range("H131").Value = 22223520
$ws.Range("I131").Value = 111111310
$ws.Range("J131").Value = 1572.5
$ws.Range("K131").Value = 333333930
$ws.Range("L131").Value = 4717.5
$ws.Range("M131").Value = -333328890
$ws.Range("N131").Value = -14797.5

$ws.Range("H134").Value = 4618.8887
$ws.Range("I134").Value = 2652.8333
$ws.Range("J134").Value = 5601.9165
$ws.Range("K134").Value = 7958.499899999999
$ws.Range("L134").Value = 16805.7495
$ws.Range("M134").Value = -2888.499899999999
$ws.Range("N134").Value = -26945.7495

$ws.Range("H137").Value = 2704.875
$ws.Range("I137").Value = 1048.4286
$ws.Range("J137").Value = 3993.2222
$ws.Range("K137").Value = 3145.2858
$ws.Range("L137").Value = 11979.6666
$ws.Range("M137").Value = 1954.7142
$ws.Range("N137").Value = -22179.6666

$ws.Range("H140").Value = 38711.168
$ws.Range("I140").Value = 49602.13
$ws.Range("K140").Value = 148806.39
$ws.Range("M140").Value = -143626.39


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2607.0833
$ws.Range("J80").Value = 3010.625
$ws.Range("L80").Value = 3010.625
$ws.Range("N80").Value = -5006.625

$ws.Range("H83").Value = 2607.0833
$ws.Range("J83").Value = 3010.625
$ws.Range("L83").Value = 15053.125
$ws.Range("N83").Value = -25037.125

$ws.Range("H109").Value = 16333.333
$ws.Range("J109").Value = 16333.333
$ws.Range("L109").Value = 16333.333
$ws.Range("N109").Value = -18413.333

$ws.Range("H132").Value = 2435.2
$ws.Range("I132").Value = 2040.7693
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 6122.3079
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -3592.3079
$ws.Range("N132").Value = -20057


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3516.4443
$ws.Range("I40").Value = 3191.3333
$ws.Range("K40").Value = 3191.3333
$ws.Range("M40").Value = -3055.3333

$ws.Range("H61").Value = 1344.9166
$ws.Range("I61").Value = 1466.875
$ws.Range("K61").Value = 1466.875
$ws.Range("M61").Value = -1264.875

$ws.Range("H113").Value = 1344.9166
$ws.Range("I113").Value = 1466.875
$ws.Range("K113").Value = 1466.875
$ws.Range("M113").Value = 703.125

$ws.Range("H132").Value = 94055.17999999999
$ws.Range("I132").Value = 3219.8
$ws.Range("J132").Value = 169751.33
$ws.Range("K132").Value = 9659.400000000001
$ws.Range("L132").Value = 509253.99
$ws.Range("M132").Value = -7129.400000000001
$ws.Range("N132").Value = -514313.99


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 29800
$ws.Range("J133").Value = 29800
$ws.Range("L133").Value = 29800
$ws.Range("N133").Value = -39920

$ws.Range("H136").Value = 1077
$ws.Range("I136").Value = 813.2857
$ws.Range("K136").Value = 2439.8571
$ws.Range("M136").Value = 110.1428999999998

